$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - SODA ½ LITRO * 12
$ws.Range("B2").Value = "SODA 500CC * 12 PET "
$ws.Range("C2").Value = 6

# Row 3 - SODA TALCA 2¼ LITRO PET* 6 U
$ws.Range("B3").Value = "SODA TALCA 2,25L PET* 6 U"
$ws.Range("C3").Value = 2

# Row 4 - SIFÓN SODA TALCA 2 LT PET*6U
$ws.Range("B4").Value = "SIFON SODA TALCA 2 LT PET*6U "
$ws.Range("C4").Value = 2

# Row 5 - TALCA COLA PET ½L DESC *12U
$ws.Range("B5").Value = "TALCA COLA PET 500CC DESC *12U "
$ws.Range("C5").Value = 0

# Row 6 - TALCA LIMA LIMON PT½L DSC*12U
$ws.Range("B6").Value = "TALCA LIMA LIMON PT 500CC DSC*12U"
$ws.Range("C6").Value = 0

# Row 7 - TALCA NARANJA PT ½L DESC*12U
$ws.Range("B7").Value = "TALCA NARANJA PT 500CC DESC*12U "
$ws.Range("C7").Value = 0

# Row 8 - TALCA POMELO PET ½L DESC *12U
$ws.Range("B8").Value = "TALCA POMELO PET 500CC DESC *12U "
$ws.Range("C8").Value = 0

# Row 9 - TALCA COLA PET 2¼L DESC *6U
$ws.Range("B9").Value = "TALCA COLA PET 2,25L DESC *6U "

# Row 10 - TALCA LIMA LIMON PT2¼LDESC*6U
$ws.Range("B10").Value = "TALCA LIMA LIMON PT 2,25L DESC*6U "

# Row 11 - TALCA NARANJA PT 2¼L DESC *6U
$ws.Range("B11").Value = "TALCA NARANJA PT 2,25L DESC *6U "

# Row 12 - TALCA POMELO PT 2¼L DESC *6U
$ws.Range("B12").Value = "TALCA POMELO PT 2,25L DESC *6U "

# Row 14 - TALCA LIMA LIMON PET 3L DESC*6 U (trailing space added)
$ws.Range("B14").Value = "TALCA LIMA LIMON PET 3L DESC*6 U "

# Row 15 - TALCA NARANJA 3L DESC *6U (trailing space added)
$ws.Range("B15").Value = "TALCA NARANJA 3L DESC *6U "
